# Weekly update to the Hortaliza / Papa (Terminal La Palmera de La Serena) sheet.
# A new price record (dated 2022-09-22 / serial 44826) is inserted as row 464,
# pushing the existing rows 464:496 down to 465:497.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 464 - shifts rows 464:496 down to 465:497
$ws.Rows.Item(464).Insert()

# Populate the newly inserted row 464 with the new record
$ws.Range("A464").Value = 8
$ws.Range("B464").Value = "Terminal La Palmera de La Serena"
$ws.Range("C464").Value = "Coquimbo"
$ws.Range("D464").Value = 44826
$ws.Range("E464").Value = 4
$ws.Range("F464").Value = 100114001
$ws.Range("G464").Value = "Papa"
$ws.Range("H464").Value = "Asterix"
$ws.Range("I464").Value = "1a (cosecha)"
$ws.Range("J464").Value = 2400
$ws.Range("K464").Value = 9000
$ws.Range("L464").Value = 10000
$ws.Range("M464").Value = 9500
$ws.Range("N464").Value = "$/saco 25 kilos"
$ws.Range("O464").Value = "Provincia de Melipilla"
$ws.Range("P464").Value = 380
$ws.Range("Q464").Value = 25
$ws.Range("R464").Value = "Hortaliza"
